$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for columns B:F, rows 2 through 18
$data = @{
  2  = @("NSE:ANANTRAJ", "NSE:ACEINTEG",  "NSE:ABFRL",    "NSE:BSOFT",       "")
  3  = @("NSE:ASMS",      "NSE:BPCL",      "NSE:PAGEIND",  "NSE:CANFINHOME",  "")
  4  = @("NSE:ASPINWALL", "NSE:COALINDIA", "",             "",                "")
  5  = @("NSE:CELEBRITY", "NSE:DPABHUSHAN","",             "",                "")
  6  = @("NSE:DALMIASUG", "NSE:DRREDDY",   "",             "",                "")
  7  = @("NSE:DCMNVL",    "NSE:GOACARBON", "",             "",                "")
  8  = @("NSE:DICIND",    "NSE:GUJGASLTD", "",             "",                "")
  9  = @("NSE:FAZE3Q",    "NSE:HEXATRADEX","",             "",                "")
  10 = @("NSE:MHRIL",     "NSE:HNDFDS",    "",             "",                "")
  11 = @("NSE:PONNIERODE","NSE:INSECTICID","",             "",                "")
  12 = @("NSE:RBL",       "NSE:JINDALPOLY","",             "",                "")
  13 = @("",              "NSE:JKPAPER",   "",             "",                "")
  14 = @("",              "NSE:LALPATHLAB","",             "",                "")
  15 = @("",              "NSE:MAANALU",   "",             "",                "")
  16 = @("",              "NSE:MARINE",    "",             "",                "")
  17 = @("",              "NSE:OIL",       "",             "",                "")
  18 = @("",              "NSE:ONWARDTEC", "",             "",                "")
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i   # B=2 .. F=6
        $cell = $ws.Cells.Item($r, $col)
        if ($vals[$i] -eq "") {
            $cell.Value = $null
        } else {
            $cell.Value = $vals[$i]
        }
    }
}

# Rows 19-22 are removed entirely (used range shrinks from A1:F22 to A1:F18)
$ws.Range("A19:F22").Clear()
